# Insert a new week's worth of records (2023-01-25, serial 44951) for
# "Black Amber" Ciruela (Primera/Segunda) above the existing rows, pushing
# the existing data rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 57:58, shifting everything below down.
$ws.Rows("57:58").Insert()

# Row 57: Black Amber / Primera
$ws.Range("A57").Value2 = 2
$ws.Range("B57").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C57").Value2 = "Coquimbo"
$ws.Range("D57").Value2 = 44951
$ws.Range("E57").Value2 = 4
$ws.Range("F57").Value2 = "Fruta"
$ws.Range("G57").Value2 = 100103
$ws.Range("H57").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I57").Value2 = 100103002
$ws.Range("J57").Value2 = "Ciruela"
$ws.Range("K57").Value2 = "Black Amber"
$ws.Range("L57").Value2 = "Primera"
$ws.Range("M57").Value2 = 16
$ws.Range("N57").Value2 = 300000
$ws.Range("O57").Value2 = 310000
$ws.Range("P57").Value2 = 305000
$ws.Range("Q57").Value2 = '$/bins (450 kilos)'
$ws.Range("R57").Value2 = "Región de O'Higgins"
$ws.Range("S57").Value2 = 678
$ws.Range("T57").Value2 = 450

# Row 58: Black Amber / Segunda
$ws.Range("A58").Value2 = 2
$ws.Range("B58").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C58").Value2 = "Coquimbo"
$ws.Range("D58").Value2 = 44951
$ws.Range("E58").Value2 = 4
$ws.Range("F58").Value2 = "Fruta"
$ws.Range("G58").Value2 = 100103
$ws.Range("H58").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I58").Value2 = 100103002
$ws.Range("J58").Value2 = "Ciruela"
$ws.Range("K58").Value2 = "Black Amber"
$ws.Range("L58").Value2 = "Segunda"
$ws.Range("M58").Value2 = 20
$ws.Range("N58").Value2 = 230000
$ws.Range("O58").Value2 = 240000
$ws.Range("P58").Value2 = 235000
$ws.Range("Q58").Value2 = '$/bins (450 kilos)'
$ws.Range("R58").Value2 = "Región de O'Higgins"
$ws.Range("S58").Value2 = 522
$ws.Range("T58").Value2 = 450
